# ADD results from server
# Update computed investment-cost results (row 2) across several year sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 7260.855929159378
$ws.Range("E2").Value = 202482.3775137067
$ws.Range("G2").Value = 56671.47998863283
$ws.Range("I2").Value = 113197.4679992575
$ws.Range("L2").Value = 406900.475060112
$ws.Range("M2").Value = 74167.63997874
$ws.Range("N2").Value = 49002.61442705191
$ws.Range("O2").Value = 48821.42540419883

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 21842.26700681437
$ws.Range("E2").Value = 119460.9857877964
$ws.Range("I2").Value = 146356.429396456
$ws.Range("L2").Value = 53298.34532628221
$ws.Range("M2").Value = 47975.70999908229
$ws.Range("N2").Value = 13499.63412324797
$ws.Range("O2").Value = 18923.19708340089

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 27543.1755456332
$ws.Range("B2").Value = 22113.21643273498
$ws.Range("E2").Value = 114655.4402706629
$ws.Range("I2").Value = 153866.0861464091
$ws.Range("M2").Value = 44638.22942194272
$ws.Range("N2").Value = 39676.88529639924
$ws.Range("O2").Value = 31311.04369977792

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 1142.580190039942
$ws.Range("O2").Value = 0

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 29588.33508286276
$ws.Range("N2").Value = 4347.543515635315
$ws.Range("O2").Value = 20429.76977394434
